$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was reported for this market/product, dated
# between the existing rows 127 and 128 (old numbering). Insert a new row
# at row 128 so every following record shifts down by one (matching the
# new dimension A1:R222), then populate the inserted row with its data.
$ws.Rows.Item(128).Insert()

$ws.Range("A128").Value = 5
$ws.Range("B128").Value = "Macroferia Regional de Talca"
$ws.Range("C128").Value = "Maule"
$ws.Range("D128").Value = 44574
$ws.Range("E128").Value = 7
$ws.Range("F128").Value = 100112006
$ws.Range("G128").Value = "Repollo"
$ws.Range("H128").Value = "Crespo record"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 3000
$ws.Range("K128").Value = 500
$ws.Range("L128").Value = 500
$ws.Range("M128").Value = 500
$ws.Range("N128").Value = "$/unidad"
$ws.Range("O128").Value = "Región del Maule"
$ws.Range("P128").Value = 500
$ws.Range("Q128").Value = 1
$ws.Range("R128").Value = "Hortaliza"
